# Apply cell value updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.665.33'
$ws.Range('E2').Value = '  -2.57%  '
$ws.Range('D3').Value = '1.984.14'
$ws.Range('E3').Value = '  -3.98%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.86'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.638'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.87%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '57.36'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +8.64%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '59.70'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0731'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('E12').Value = '  -4.53%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.923'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.61%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.06'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').Value = '2.276.86'
$ws.Range('E15').Value = '  -3.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.25'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('D17').Value = '1.982.05'
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '17.23'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +5.02%  '
$ws.Range('D19').Value = '35.552.95'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '233.25'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.06'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.00%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.29'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.77%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.14'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.51'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.45'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.55%  '
$ws.Range('E30').Value = '  -3.83%  '
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.80'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -5.38%  '
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0899'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +10.09%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.37'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.27'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.78%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.80'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.92'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('E40').Value = '  -4.88%  '
$ws.Range('E41').Value = '  -2.25%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0210'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('E43').Value = '  -4.07%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0890'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.37%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.98'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').Value = '1.377.15'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '15.45'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.89'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.25'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '45.79'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.66%  '
